# Update "Producción industrial INE (base 2014=100) desestacionalizada"
# (column C) values for several years, and one value in column E (2019),
# as per the source data refresh ("Actualización desde MV -datos-").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = 88.16
$ws.Range("C12").Value = 90.14
$ws.Range("C13").Value = 93.2
$ws.Range("C14").Value = 96.41
$ws.Range("C15").Value = 99.94
$ws.Range("C17").Value = 100.53
$ws.Range("C19").Value = 98.56999999999999
$ws.Range("C20").Value = 102.17
$ws.Range("C21").Value = 101.47
$ws.Range("E21").Value = 98.43000000000001
